$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha (D) and Volumen (J) values between rows 3 and 4
$ws.Range("D3").Value = 44414
$ws.Range("J3").Value = 500

$ws.Range("D4").Value = 44379
$ws.Range("J4").Value = 240
